$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2,3,5,6,7,8) get their contents cyclically rotated between
# rows, while row 4 stays untouched. Capture the "before" values for the
# columns that change (D, I, J, K, L, M, N, P, Q) first, then write them
# back out according to the permutation described by the diff:
#   new row2 <- old row6
#   new row3 <- old row5
#   new row5 <- old row3
#   new row6 <- old row8
#   new row7 <- old row2
#   new row8 <- old row7

$cols = @("D","I","J","K","L","M","N","P","Q")

$rows = @(2, 3, 4, 5, 6, 7, 8)
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

$mapping = @{ 2 = 6; 3 = 5; 5 = 3; 6 = 8; 7 = 2; 8 = 7 }

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $snapshot[$srcRow][$c]
    }
}
